# Updates cryptos list values (Price / Volume(1h)) per the Sun Sep 10
# 22:21:34 UTC 2023 GitHub Actions refresh, plus a ranking swap between
# TrustWalletToken and RocketPoolETH (rows 42/43).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "26.031.08"
$ws.Range("E2").Value = "  -0.06%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.630.38"
$ws.Range("E3").Value = "  -0.90%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.03%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'214.36"
$ws.Range("E5").Value = "  -0.60%  "

# Row 6 - XRP (price unchanged)
$ws.Range("E6").Value = "  -0.70%  "

# Row 7 - USDC
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  -0.09%  "

# Row 8 - Cardano (price unchanged)
$ws.Range("E8").Value = "  -1.83%  "

# Row 9 - Dogecoin
$ws.Range("D9").Value = "'0.0620"
$ws.Range("E9").Value = "  -3.01%  "

# Row 10 - Solana
$ws.Range("D10").Value = "'18.62"
$ws.Range("E10").Value = "  -5.17%  "

# Row 11 - TRON
$ws.Range("D11").Value = "'0.0789"
$ws.Range("E11").Value = "  -0.96%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.856.90"
$ws.Range("E12").Value = "  -0.91%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.640.14"
$ws.Range("E13").Value = "  +1.07%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "'4.19"
$ws.Range("E14").Value = "  -1.86%  "

# Row 15 - Polygon
$ws.Range("D15").Value = "'0.529"
$ws.Range("E15").Value = "  -2.78%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "26.032.07"
$ws.Range("E16").Value = "  -0.07%  "

# Row 17 - ShibaInu
$ws.Range("D17").Value = "0.0₃0742"
$ws.Range("E17").Value = "  -2.50%  "

# Row 18 - Litecoin
$ws.Range("D18").Value = "'61.62"
$ws.Range("E18").Value = "  -2.96%  "

# Row 19 - Dai (price unchanged)
$ws.Range("E19").Value = "  -0.05%  "

# Row 20 - BitcoinCash (volume unchanged)
$ws.Range("D20").Value = "'193.15"

# Row 21 - Uniswap (price unchanged)
$ws.Range("E21").Value = "  -2.19%  "

# Row 22 - Avalanche
$ws.Range("D22").Value = "'9.56"
$ws.Range("E22").Value = "  -3.61%  "

# Row 23 - Chainlink
$ws.Range("D23").Value = "'6.08"
$ws.Range("E23").Value = "  -2.09%  "

# Row 24 - Stellar (price unchanged)
$ws.Range("E24").Value = "  +1.54%  "

# Row 25 - Monero
$ws.Range("D25").Value = "'144.24"
$ws.Range("E25").Value = "  +0.09%  "

# Row 26 - BinanceUSD (price unchanged)
$ws.Range("E26").Value = "  -0.04%  "

# Row 27 - Toncoin (price unchanged)
$ws.Range("E27").Value = "  -3.65%  "

# Row 28 - Cosmos
$ws.Range("D28").Value = "'6.74"
$ws.Range("E28").Value = "  -2.15%  "

# Row 29 - EthereumClassic
$ws.Range("D29").Value = "'15.30"
$ws.Range("E29").Value = "  -1.40%  "

# Row 30 - PancakeSwap (price unchanged)
$ws.Range("E30").Value = "  -0.64%  "

# Row 31 - Hedera (price unchanged)
$ws.Range("E31").Value = "  -2.45%  "

# Row 32 - Filecoin
$ws.Range("D32").Value = "'3.14"
$ws.Range("E32").Value = "  -4.05%  "

# Row 33 - InternetComputer(DFINITY) (price unchanged)
$ws.Range("E33").Value = "  -4.74%  "

# Row 34 - LidoDAOToken (price unchanged)
$ws.Range("E34").Value = "  -3.06%  "

# Row 35 - HuobiToken (price unchanged)
$ws.Range("E35").Value = "  -1.98%  "

# Row 36 - Maker
$ws.Range("D36").Value = "1.124.32"
$ws.Range("E36").Value = "  -0.73%  "

# Row 37 - ARBITRUM
$ws.Range("D37").Value = "'0.854"
$ws.Range("E37").Value = "  -5.65%  "

# Row 38 - MXToken
$ws.Range("D38").Value = "'2.43"
$ws.Range("E38").Value = "  -1.24%  "

# Row 39 - ImmutableX
$ws.Range("D39").Value = "'0.521"
$ws.Range("E39").Value = "  -3.34%  "

# Row 40 - VeChain
$ws.Range("D40").Value = "'0.0154"
$ws.Range("E40").Value = "  -2.32%  "

# Row 41 - Quant (price unchanged)
$ws.Range("E41").Value = "  -0.61%  "

# Rows 42/43 - ranking swap: RocketPoolETH now ranks above TrustWalletToken
$ws.Range("B42").Value = "RocketPoolETH"
$ws.Range("C42").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D42").Value = "1.766.49"
$ws.Range("E42").Value = "  -0.95%  "

$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "'0.763"
$ws.Range("E43").Value = "  -4.13%  "

# Row 44 - FraxShare
$ws.Range("D44").Value = "'5.14"
$ws.Range("E44").Value = "  -5.75%  "

# Row 45 - BabyDogeCoin (price unchanged)
$ws.Range("E45").Value = "  -1.84%  "

# Row 46 - Aave
$ws.Range("D46").Value = "'54.66"
$ws.Range("E46").Value = "  -3.32%  "

# Row 47 - Cronos
$ws.Range("D47").Value = "'0.0525"
$ws.Range("E47").Value = "  +0.56%  "

# Row 48 - RenderToken (price unchanged)
$ws.Range("E48").Value = "  -0.38%  "

# Row 49 - Mantle
$ws.Range("D49").Value = "'0.413"
$ws.Range("E49").Value = "  -0.35%  "

# Row 50 - EnergySwap (price unchanged)
$ws.Range("E50").Value = "  -3.17%  "

# Row 51 - USDD
$ws.Range("D51").Value = "'1.00"
$ws.Range("E51").Value = "  +0.03%  "
